# Rename "Sheet1" to "F-SW-SD-06"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "F-SW-SD-06"

# Update the workbook-level Print_Area defined name so it points at the
# renamed sheet instead of the old "Sheet1" reference.
foreach ($dn in $wb.Names) {
    if ($dn.Name -eq "F-SW-SD-06!Print_Area") {
        $dn.RefersTo = "='F-SW-SD-06'!`$A`$1:`$G`$31"
    }
}
